$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update bus voltage magnitude results (p.u.) for the 380 kV case.
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.04506430032093
$ws.Range("D2").Value = 1.051698405279343
$ws.Range("E2").Value = 1.052795448883082
$ws.Range("F2").Value = 1.063451034760574
$ws.Range("I2").Value = 1.045767134120241
$ws.Range("J2").Value = 1.050126033778945
$ws.Range("K2").Value = 1.054449185893118
$ws.Range("L2").Value = 1.055543193097044
$ws.Range("M2").Value = 1.066169639263911
$ws.Range("N2").Value = 1.05161733319805
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.045925137196454
$ws.Range("D3").Value = 1.052374087308295
$ws.Range("E3").Value = 1.053558308311489
$ws.Range("F3").Value = 1.064288539036854
$ws.Range("I3").Value = 1.04599613584192
$ws.Range("J3").Value = 1.050634818944301
$ws.Range("K3").Value = 1.054937918421713
$ws.Range("L3").Value = 1.056119097755566
$ws.Range("M3").Value = 1.066822100355134
$ws.Range("N3").Value = 1.0521268408967
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04648273790602
$ws.Range("D4").Value = 1.052811775128941
$ws.Range("E4").Value = 1.054052829295803
$ws.Range("F4").Value = 1.064831462586124
$ws.Range("I4").Value = 1.046143257411248
$ws.Range("J4").Value = 1.050963946923528
$ws.Range("K4").Value = 1.055253955041966
$ws.Range("L4").Value = 1.056491981453951
$ws.Range("M4").Value = 1.067244637355262
$ws.Range("N4").Value = 1.052456436275409
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.046717290916649
$ws.Range("D5").Value = 1.052995891382413
$ws.Range("E5").Value = 1.054260939553609
$ws.Range("F5").Value = 1.065059945679874
$ws.Range("I5").Value = 1.046204853261524
$ws.Range("J5").Value = 1.05110228949856
$ws.Range("K5").Value = 1.055386766210052
$ws.Range("L5").Value = 1.056648796567911
$ws.Range("M5").Value = 1.067422354356361
$ws.Range("N5").Value = 1.052594975312765
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.046756681426368
$ws.Range("D6").Value = 1.053026811845578
$ws.Range("E6").Value = 1.054295894656768
$ws.Range("F6").Value = 1.065098322899427
$ws.Range("I6").Value = 1.046215180554415
$ws.Range("J6").Value = 1.051125516450095
$ws.Range("K6").Value = 1.055409062782723
$ws.Range("L6").Value = 1.056675129688124
$ws.Range("M6").Value = 1.067452198611154
$ws.Range("N6").Value = 1.052618235249236
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04648587147585
$ws.Range("D7").Value = 1.052814234860294
$ws.Range("E7").Value = 1.054055609238082
$ws.Range("F7").Value = 1.0648345146555
$ws.Range("I7").Value = 1.046144081457289
$ws.Range("J7").Value = 1.050965795555076
$ws.Range("K7").Value = 1.055255729871865
$ws.Range("L7").Value = 1.05649407661161
$ws.Range("M7").Value = 1.067247011696574
$ws.Range("N7").Value = 1.052458287532226
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045355102777114
$ws.Range("D8").Value = 1.051926655498481
$ws.Range("E8").Value = 1.053053073377221
$ws.Range("F8").Value = 1.063733864572442
$ws.Range("I8").Value = 1.045844744876727
$ws.Range("J8").Value = 1.050297997987951
$ws.Range("K8").Value = 1.054614397207575
$ws.Range("L8").Value = 1.05573777328538
$ws.Range("M8").Value = 1.066390067899581
$ws.Range("N8").Value = 1.051789541615958
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.043367069041824
$ws.Range("D9").Value = 1.050366352950183
$ws.Range("E9").Value = 1.051293449323716
$ws.Range("F9").Value = 1.061802139398889
$ws.Range("I9").Value = 1.045309209241289
$ws.Range("J9").Value = 1.049120619893773
$ws.Range("K9").Value = 1.053482767197605
$ws.Range("L9").Value = 1.054406931924681
$ws.Range("M9").Value = 1.064882782523879
$ws.Range("N9").Value = 1.05061049150984
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042044849914803
$ws.Range("D10").Value = 1.049328760406604
$ws.Range("E10").Value = 1.050125162072928
$ws.Range("F10").Value = 1.0605196547005
$ws.Range("I10").Value = 1.044946810026286
$ws.Range("J10").Value = 1.048335350879623
$ws.Range("K10").Value = 1.052727403195824
$ws.Range("L10").Value = 1.053521040779255
$ws.Range("M10").Value = 1.063879875654256
$ws.Range("N10").Value = 1.049824107323618
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041473078170189
$ws.Range("D11").Value = 1.048880112224154
$ws.Range("E11").Value = 1.049620440947355
$ws.Range("F11").Value = 1.059965613437654
$ws.Range("I11").Value = 1.044788622231747
$ws.Range("J11").Value = 1.047995253551251
$ws.Range("K11").Value = 1.05240011474165
$ws.Range("L11").Value = 1.053137774963548
$ws.Range("M11").Value = 1.063446088454769
$ws.Range("N11").Value = 1.049483527018029
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041260812122976
$ws.Range("D12").Value = 1.048713561676464
$ws.Range("E12").Value = 1.049433140409721
$ws.Range("F12").Value = 1.059760012559265
$ws.Range("I12").Value = 1.044729674747689
$ws.Range("J12").Value = 1.047868916917215
$ws.Range("K12").Value = 1.052278514871461
$ws.Range("L12").Value = 1.052995464123998
$ws.Range("M12").Value = 1.063285034015106
$ws.Range("N12").Value = 1.049357010971482
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041306338661549
$ws.Range("D13").Value = 1.048749282899923
$ws.Range("E13").Value = 1.04947330904585
$ws.Range("F13").Value = 1.059804105803066
$ws.Range("I13").Value = 1.044742327745418
$ws.Range("J13").Value = 1.047896016956447
$ws.Range("K13").Value = 1.052304599809433
$ws.Range("L13").Value = 1.053025987944504
$ws.Range("M13").Value = 1.0633195773899
$ws.Range("N13").Value = 1.049384149495878
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041455529828502
$ws.Range("D14").Value = 1.048866343107479
$ws.Range("E14").Value = 1.04960495503853
$ws.Range("F14").Value = 1.059948614421513
$ws.Range("I14").Value = 1.044783753475368
$ws.Range("J14").Value = 1.04798481070891
$ws.Range("K14").Value = 1.052390063870307
$ws.Range("L14").Value = 1.053126010453423
$ws.Range("M14").Value = 1.063432774137055
$ws.Range("N14").Value = 1.049473069345655
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.041547466774116
$ws.Range("D15").Value = 1.048938480719507
$ws.Range("E15").Value = 1.049686089771576
$ws.Range("F15").Value = 1.060037676810643
$ws.Range("I15").Value = 1.044809252148697
$ws.Range("J15").Value = 1.048039518270352
$ws.Range("K15").Value = 1.052442717114394
$ws.Range("L15").Value = 1.053187644443053
$ws.Range("M15").Value = 1.06350252816876
$ws.Range("N15").Value = 1.049527854598109
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042082812624202
$ws.Range("D16").Value = 1.049358549265298
$ws.Range("E16").Value = 1.0501586832427
$ws.Range("F16").Value = 1.060556451793621
$ws.Range("I16").Value = 1.044957281812011
$ws.Range("J16").Value = 1.048357920639077
$ws.Range("K16").Value = 1.052749119922202
$ws.Range("L16").Value = 1.053546483963368
$ws.Range("M16").Value = 1.063908674892569
$ws.Range("N16").Value = 1.049846709134719
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.042418824945594
$ws.Range("D17").Value = 1.049622218758874
$ws.Range("E17").Value = 1.050455439157945
$ws.Range("F17").Value = 1.060882210525456
$ws.Range("I17").Value = 1.045049798376481
$ws.Range("J17").Value = 1.048557628012843
$ws.Range("K17").Value = 1.0529412627831
$ws.Range("L17").Value = 1.053771664118149
$ws.Range("M17").Value = 1.064163569139759
$ws.Range("N17").Value = 1.050046700115863
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.042614888264641
$ws.Range("D18").Value = 1.049776073982096
$ws.Range("E18").Value = 1.050628643176459
$ws.Range("F18").Value = 1.061072343749858
$ws.Range("I18").Value = 1.045103639467182
$ws.Range("J18").Value = 1.048674107041822
$ws.Range("K18").Value = 1.053053316030392
$ws.Range("L18").Value = 1.053903039723105
$ws.Range("M18").Value = 1.064312290662629
$ws.Range("N18").Value = 1.050163344558422
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.042681753126714
$ws.Range("D19").Value = 1.049828544975383
$ws.Range("E19").Value = 1.050687720095511
$ws.Range("F19").Value = 1.061137195182778
$ws.Range("I19").Value = 1.045121977126971
$ws.Range("J19").Value = 1.048713822170962
$ws.Range("K19").Value = 1.053091519786318
$ws.Range("L19").Value = 1.053947840764645
$ws.Range("M19").Value = 1.064363008615505
$ws.Range("N19").Value = 1.050203116087602
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.042382766439638
$ws.Range("D20").Value = 1.049593923167445
$ws.Range("E20").Value = 1.05042358852971
$ws.Range("F20").Value = 1.060847246890912
$ws.Range("I20").Value = 1.045039884864503
$ws.Range("J20").Value = 1.048536202006509
$ws.Range("K20").Value = 1.052920649776518
$ws.Range("L20").Value = 1.053747501115892
$ws.Range("M20").Value = 1.06413621662753
$ws.Range("N20").Value = 1.050025243682141
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041411593549026
$ws.Range("D21").Value = 1.048831869092045
$ws.Range("E21").Value = 1.049566183719275
$ws.Range("F21").Value = 1.059906054844429
$ws.Range("I21").Value = 1.044771559855013
$ws.Range("J21").Value = 1.047958663409637
$ws.Range("K21").Value = 1.052364897662215
$ws.Range("L21").Value = 1.053096554894041
$ws.Range("M21").Value = 1.063399438476105
$ws.Range("N21").Value = 1.049446884914217
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040801647272034
$ws.Range("D22").Value = 1.048353299640929
$ws.Range("E22").Value = 1.049028114780141
$ws.Range("F22").Value = 1.059315417575358
$ws.Range("I22").Value = 1.044601757311042
$ws.Range("J22").Value = 1.047595488543407
$ws.Range("K22").Value = 1.052015299213059
$ws.Range("L22").Value = 1.052687575813047
$ws.Range("M22").Value = 1.062936622512228
$ws.Range("N22").Value = 1.049083194298021
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.041124927424973
$ws.Range("D23").Value = 1.0486069442904
$ws.Range("E23").Value = 1.049313258433831
$ws.Range("F23").Value = 1.059628418040265
$ws.Range("I23").Value = 1.044691876472564
$ws.Range("J23").Value = 1.047788019050715
$ws.Range("K23").Value = 1.052200644118646
$ws.Range("L23").Value = 1.052904354764175
$ws.Range("M23").Value = 1.063181929121442
$ws.Range("N23").Value = 1.049275998220733
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.042399059500929
$ws.Range("D24").Value = 1.049606708538156
$ws.Range("E24").Value = 1.050437980113907
$ws.Range("F24").Value = 1.060863045070504
$ws.Range("I24").Value = 1.045044364731685
$ws.Range("J24").Value = 1.048545883517982
$ws.Range("K24").Value = 1.052929963970208
$ws.Range("L24").Value = 1.053758419238353
$ws.Range("M24").Value = 1.064148575908433
$ws.Range("N24").Value = 1.050034938942471
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.043880477460844
$ws.Range("D25").Value = 1.050769276300649
$ws.Range("E25").Value = 1.051747516970345
$ws.Range("F25").Value = 1.062300605193824
$ws.Range("I25").Value = 1.045448609013341
$ws.Range("J25").Value = 1.049425067630118
$ws.Range("K25").Value = 1.053775492769696
$ws.Range("L25").Value = 1.054750756962329
$ws.Range("M25").Value = 1.065272115167433
$ws.Range("N25").Value = 1.050915371596891
